$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C holds a "Förändrad" (last changed) date for each record, stored
# as serial date 45171 (2023-09-02). Bump it by one day to 45172
# (2023-09-03) for every data row (rows 2 through 321).
$ws.Range("C2:C321").Value = 45172
